$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (dates) from serial 45177 to 45178 for rows 2 through 39
for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value = 45178
    }
}
